# Beluga_data.xlsx - "Strandings" sheet: add a second carcass-composition
# model that is fit to Newborn/Older removals without using the age-based
# weighting used by the first model (0.8/0.2), instead using (0.83/0.17)
# for newborns (column D, from column B) and (0.72/0.28) for "older"
# (column E, from column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strandings")

# --- Row 3 (first data row with a formula) ---------------------------------
$ws.Range("D3").Formula = "=ROUND(B3*0.83,0)+ROUND(B2*0.17,0)"
$ws.Range("E3").Formula = "=ROUND(C3*0.72,0)+ROUND(C2*0.28,0)"

# --- Rows 4:32 - fill the updated formulas down ----------------------------
$ws.Range("D4:D32").Formula = "=ROUND(B4*0.83,0)+ROUND(B3*0.17,0)"
$ws.Range("E4:E32").Formula = "=ROUND(C4*0.72,0)+ROUND(C3*0.28,0)"

# --- Rows 30 & 32 previously held hard-coded carcass-age overrides; the new
#     model drops the age adjustment there entirely and just carries the
#     newborn removals straight down from column B. ------------------------
$ws.Range("D30").Formula = "=B30"
$ws.Range("D32").Formula = "=B32"

# --- Selection state, to mirror the saved workbook's cursor position -------
$ws.Range("A33:E34").Select()
$excel.ActiveWindow.RangeSelection.Item(1,1).Activate() | Out-Null
